$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumToText($cellRef, $donorRef) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = "0"
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

function Set-TextToNum($cellRef, $donorRef, $value) {
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($cellRef).Value = $value
}

# --- Shared-string text edits (Volume number, report date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "48"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "11/28/2022"
$c9 = $ws.Range("C9")
$c9.Characters(48, 10).Text = "12/4/2022"

# --- Cell type conversions (number <-> text placeholder) ---
Set-NumToText "C15" "D15"
Set-NumToText "C26" "D26"
Set-NumToText "F30" "C30"

Set-TextToNum "C22" "D16" 2
Set-TextToNum "D27" "D16" 2
Set-TextToNum "E27" "E16" 0
Set-TextToNum "D30" "D16" 3
Set-TextToNum "E30" "E16" -100

# --- Plain numeric value updates ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -52.380952380952
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 126
$ws.Range("K16").Value = 1.587301587301
$ws.Range("L16").Value = -5.882352941176
$ws.Range("M16").Value = 52.380952380952
$ws.Range("N16").Value = -84.039900249376
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 114
$ws.Range("K17").Value = 9.649122807017
$ws.Range("L17").Value = 64.473684210526
$ws.Range("M17").Value = 92.307692307692
$ws.Range("N17").Value = -27.325581395348
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 18
$ws.Range("H18").Value = -28
$ws.Range("I18").Value = 236
$ws.Range("J18").Value = 170
$ws.Range("K18").Value = 38.823529411764
$ws.Range("L18").Value = 3.056768558951
$ws.Range("M18").Value = 47.5
$ws.Range("N18").Value = -71.904761904761
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 92
$ws.Range("G19").Value = 100
$ws.Range("H19").Value = -8
$ws.Range("I19").Value = 1133
$ws.Range("J19").Value = 781
$ws.Range("K19").Value = 45.070422535211
$ws.Range("L19").Value = 82.154340836012
$ws.Range("M19").Value = 13.186813186813
$ws.Range("N19").Value = -68.718939812258
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 700
$ws.Range("I20").Value = 70
$ws.Range("K20").Value = 7.692307692307
$ws.Range("L20").Value = 125.806451612903
$ws.Range("M20").Value = 89.189189189189
$ws.Range("N20").Value = -91.764705882352
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -11.764705882352
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 158
$ws.Range("H21").Value = -10.759493670886
$ws.Range("I21").Value = 1713
$ws.Range("J21").Value = 1271
$ws.Range("K21").Value = 34.775767112509
$ws.Range("L21").Value = 53.770197486535
$ws.Range("M21").Value = 26.60753880266
$ws.Range("N21").Value = -72.844007609384
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = -40
$ws.Range("I22").Value = 89
$ws.Range("J22").Value = 77
$ws.Range("K22").Value = 15.584415584415
$ws.Range("L22").Value = 11.25
$ws.Range("M22").Value = 43.548387096774
$ws.Range("C24").Value = 91
$ws.Range("D24").Value = 55
$ws.Range("E24").Value = 65.454545454545
$ws.Range("F24").Value = 363
$ws.Range("G24").Value = 229
$ws.Range("H24").Value = 58.515283842794
$ws.Range("I24").Value = 3736
$ws.Range("J24").Value = 2070
$ws.Range("K24").Value = 80.483091787439
$ws.Range("L24").Value = 137.96178343949
$ws.Range("M24").Value = 135.116425424795
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -6.666666666666
$ws.Range("I25").Value = 346
$ws.Range("J25").Value = 336
$ws.Range("K25").Value = 2.97619047619
$ws.Range("L25").Value = 42.38683127572
$ws.Range("M25").Value = 42.97520661157
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 110
$ws.Range("J27").Value = 77
$ws.Range("K27").Value = 42.857142857142
$ws.Range("L27").Value = 66.666666666666
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 14
$ws.Range("K30").Value = 7.142857142857
